$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Leap2","Ghsr","ECs",3,1,2.537013666666667,7.611041,0.2842954840254888,0.2842954840254888,2,0.6666666666666666,0.4759383333333333,1.427815,0.1284588056245616,0.1284588056245616,1.207462056157222,10.867158505415,0.03652025832237094,0.03652025832237092),
  @("ECs","Leap2","Ghsr","FAPs",3,1,2.537013666666667,7.611041,0.2842954840254888,0.2842954840254888,3,1,2.371854333333333,7.115563,0.6401786816403544,0.6401786816403543,6.017426859009222,54.156841731083,0.1819999081597439,0.1819999081597438),
  @("ECs","Leap2","Ghsr","M2",3,1,2.537013666666667,7.611041,0.2842954840254888,0.2842954840254888,1,0.3333333333333333,0.055189,0.165567,0.01489586470995318,0.01489586470995318,0.1400152472496667,1.260137225247,0.004234827067694335,0.004234827067694335),
  @("ECs","Leap2","Ghsr","sCs",3,1,2.537013666666667,7.611041,0.2842954840254888,0.2842954840254888,3,1,0.8020063333333334,2.406019,0.2164666480251308,0.2164666480251308,2.034701028419889,18.312309255779,0.06154049047567969,0.06154049047567967),
  @("FAPs","Leap2","Ghsr","ECs",3,1,2.670935,8.012805,0.2993025889463553,0.2993025889463553,2,0.6666666666666666,0.4759383333333333,1.427815,0.1284588056245616,0.1284588056245616,1.271200352341667,11.440803171075,0.03844805309638793,0.03844805309638791),
  @("FAPs","Leap2","Ghsr","FAPs",3,1,2.670935,8.012805,0.2993025889463553,0.2993025889463553,3,1,2.371854333333333,7.115563,0.6401786816403544,0.6401786816403543,6.335068753801667,57.015618784215,0.1916071368032227,0.1916071368032226),
  @("FAPs","Leap2","Ghsr","M2",3,1,2.670935,8.012805,0.2993025889463553,0.2993025889463553,1,0.3333333333333333,0.055189,0.165567,0.01489586470995318,0.01489586470995318,0.147406231715,1.326656085435,0.004458370872283636,0.004458370872283635),
  @("FAPs","Leap2","Ghsr","sCs",3,1,2.670935,8.012805,0.2993025889463553,0.2993025889463553,3,1,0.8020063333333334,2.406019,0.2164666480251308,0.2164666480251308,2.142106785921667,19.278961073295,0.06478902817446111,0.06478902817446108),
  @("M2","Leap2","Ghsr","ECs",3,1,1.585151666666667,4.755455,0.1776306790341073,0.1776306790341073,2,0.6666666666666666,0.4759383333333333,1.427815,0.1284588056245616,0.1284588056245616,0.7544344423138889,6.789909980825,0.02281822487100128,0.02281822487100128),
  @("M2","Leap2","Ghsr","FAPs",3,1,1.585151666666667,4.755455,0.1776306790341073,0.1776306790341073,3,1,2.371854333333333,7.115563,0.6401786816403544,0.6401786816403543,3.759748849573889,33.83773964616501,0.1137153739229358,0.1137153739229357),
  @("M2","Leap2","Ghsr","M2",3,1,1.585151666666667,4.755455,0.1776306790341073,0.1776306790341073,1,0.3333333333333333,0.055189,0.165567,0.01489586470995318,0.01489586470995318,0.08748293533166666,0.787346417985,0.002645962563229178,0.002645962563229178),
  @("M2","Leap2","Ghsr","sCs",3,1,1.585151666666667,4.755455,0.1776306790341073,0.1776306790341073,3,1,0.8020063333333334,2.406019,0.2164666480251308,0.2164666480251308,1.271301675960556,11.441715083645,0.03845111767694109,0.03845111767694108),
  @("sCs","Leap2","Ghsr","ECs",3,1,2.130761666666667,6.392285,0.2387712479940486,0.2387712479940486,2,0.6666666666666666,0.4759383333333333,1.427815,0.1284588056245616,0.1284588056245616,1.014111156363889,9.127000407274998,0.0306722693348015,0.03067226933480148),
  @("sCs","Leap2","Ghsr","FAPs",3,1,2.130761666666667,6.392285,0.2387712479940486,0.2387712479940486,3,1,2.371854333333333,7.115563,0.6401786816403544,0.6401786816403543,5.053856292383889,45.484706631455,0.1528562627544522,0.1528562627544521),
  @("sCs","Leap2","Ghsr","M2",3,1,2.130761666666667,6.392285,0.2387712479940486,0.2387712479940486,1,0.3333333333333333,0.055189,0.165567,0.01489586470995318,0.01489586470995318,0.1175946056216667,1.058351450595,0.003556704206746027,0.003556704206746026),
  @("sCs","Leap2","Ghsr","sCs",3,1,2.130761666666667,6.392285,0.2387712479940486,0.2387712479940486,3,1,0.8020063333333334,2.406019,0.2164666480251308,0.2164666480251308,1.708884351490556,15.379959163415,0.05168601169804895,0.05168601169804894)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $rowData = $data[$r]
  for ($c = 0; $c -lt $rowData.Length; $c++) {
    $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
  }
}
